$s99 = 'http://localhost:8080/invitation/sendBoardInvitation'
$s100 = 'To send Board invitations'
$s101 = @'
{
    "emailIdList": [
        "santhosh.tndr@gmail.com",
        "rajeswar061965@gmail.com"
    ],
    "emailBody": "You are invited to Borad",
    "emailSubject": "Borad Invitation",
    "createdBy": 1
}
'@
$s102 = @'
[
    {
        "email": "santhosh.tndr@gmail.com",
        "status": "SUCCESS",
        "message": null
    },
    {
        "email": "rajeswar061965@gmail.com",
        "status": "SUCCESS",
        "message": null
    }
]
'@
$s103 = 'http://localhost:8080/invitation/list'
$s104 = 'To get all Invitations'
$s105 = @'
[
    {
        "inviteeName": null,
        "email": "santhosh.tndr@gmail.com",
        "subject": "Borad Invitation",
        "message": "You are invited to Borad",
        "status": "SUCCESS",
        "statusMsg": null,
        "createdBy": 1,
        "createdDate": "2020-11-29T12:50:12.000+00:00",
        "userName": "NotificationBoard"
    },
    {
        "inviteeName": null,
        "email": "rajeswar061965@gmail.com",
        "subject": "Borad Invitation",
        "message": "You are invited to Borad",
        "status": "SUCCESS",
        "statusMsg": null,
        "createdBy": 1,
        "createdDate": "2020-11-29T12:50:13.000+00:00",
        "userName": "NotificationBoard"
    }
]
'@
$s106 = 'http://localhost:8080/user/register'
$s107 = @'
To register User (for user type 
member groupName is required
'@
$s108 = @'
{
    "userName": "NotificationBoard",
    "password": "NB@2020",
    "email": "notificationboard1tts@gmail.com",
    "alternateEmail": "notificationboardalt1@gmail.com",
    "contactNumber": "9874563210",
    "userType": "Member",
    "groupName": "UGroup1"
}
'@
$s109 = @'
{
    "message": "Created Successfully",
    "results": {
        "userId": 21,
        "userName": "NotificationBoard",
        "password": "NB@2020",
        "email": "notificationboard1tts@gmail.com",
        "alternateEmail": "notificationboardalt1@gmail.com",
        "contactNumber": "9874563210",
        "userType": "Member",
        "groupName": "UGroup1",
        "createdDate": null,
        "updatedDate": null,
        "permissions": null,
        "isActive": null,
        "authorities": null
    }
}
'@
$s110 = 'http://localhost:8080/user/{useremail}'
$s111 = 'To get user by email'
$s112 = 'http://localhost:8080/user/update'
$s113 = 'To update user'
$s114 = @'
{
    "userId": 21,
    "userName": "NotificationBoard",
    "password": "NB@2020",
    "email": "notificationboard1ttsa@gmail.com",
    "alternateEmail": "notificationboardalt1@gmail.com",
    "contactNumber": "9874563210",
    "userType": "Member",
    "groupName": "UGroup1"
}
'@
$s115 = @'
{
    "message": "Updated Successfully",
    "results": {
        "userId": 21,
        "userName": "NotificationBoard",
        "password": "NB@2020",
        "email": "notificationboard1ttsa@gmail.com",
        "alternateEmail": "notificationboardalt1@gmail.com",
        "contactNumber": "9874563210",
        "userType": "Member",
        "groupName": "UGroup1",
        "createdDate": null,
        "updatedDate": null,
        "permissions": null,
        "isActive": null,
        "authorities": null
    }
}
'@
$s116 = 'http://localhost:8080/user/delete'
$s117 = @'
{
    "email": "notificationboard1tts@gmail.com"
}
'@
$s118 = 'To delete user'
$s119 = @'
{
    "message": "Deleted Successfully",
    "results": true
}
'@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Set cell values in the exact order needed to reproduce shared-string order ----

# Row 13
$ws.Range("A13").Value = $s99
$ws.Range("C13").Value = $s100
$ws.Range("D13").Value = $s101
$ws.Range("E13").Value = $s102
$ws.Range("B13").Value = "POST"

# Row 14
$ws.Range("A14").Value = $s103
$ws.Range("C14").Value = $s104
$ws.Range("E14").Value = $s105
$ws.Range("B14").Value = "GET"

# Row 15
$ws.Range("A15").Value = $s106
$ws.Range("C15").Value = $s107
$ws.Range("D15").Value = $s108
$ws.Range("E15").Value = $s109
$ws.Range("B15").Value = "POST"

# Row 16
$ws.Range("A16").Value = $s110
$ws.Range("C16").Value = $s111
$ws.Range("B16").Value = "GET"

# Row 17
$ws.Range("A17").Value = $s112
$ws.Range("C17").Value = $s113
$ws.Range("D17").Value = $s114
$ws.Range("E17").Value = $s115
$ws.Range("B17").Value = "POST"

# Row 18 (note: request body D18 entered before purpose C18, matching original authoring order)
$ws.Range("A18").Value = $s116
$ws.Range("D18").Value = $s117
$ws.Range("C18").Value = $s118
$ws.Range("E18").Value = $s119
$ws.Range("B18").Value = "POST"

# ---- Formatting: vertical-center alignment across all new data rows ----
$ws.Range("A13:E18").VerticalAlignment = -4108

# ---- Wrap text for the long JSON / multi-line cells ----
$ws.Range("D13:E13").WrapText = $true
$ws.Range("D14:E14").WrapText = $true
$ws.Range("C15:E15").WrapText = $true
$ws.Range("D17:E17").WrapText = $true
$ws.Range("D18:E18").WrapText = $true

# ---- Hyperlinks for the new URL cells ----
$ws.Hyperlinks.Add($ws.Range("A13"), $s99)
$ws.Hyperlinks.Add($ws.Range("A14"), $s103)
$ws.Hyperlinks.Add($ws.Range("A15"), $s106)
$ws.Hyperlinks.Add($ws.Range("A16"), $s110)
$ws.Hyperlinks.Add($ws.Range("A17"), $s112)
$ws.Hyperlinks.Add($ws.Range("A18"), $s116)
$ws.Range("A13:A18").VerticalAlignment = -4108

# ---- Row heights ----
$ws.Rows.Item(13).RowHeight = 180
$ws.Rows.Item(14).RowHeight = 360
$ws.Rows.Item(15).RowHeight = 270
$ws.Rows.Item(17).RowHeight = 270
$ws.Rows.Item(18).RowHeight = 60

# ---- Rows 19-21: blank formatted cells in column C ----
$ws.Range("C19").VerticalAlignment = -4108
$ws.Range("C20").VerticalAlignment = -4108
$ws.Range("C21").VerticalAlignment = -4108

# ---- Column widths ----
$ws.Columns.Item(4).ColumnWidth = 54.71
$ws.Columns.Item(5).ColumnWidth = 81.17

# ---- View: selection / scroll ----
$ws.Range("C19").Select()
try {
  $excel.ActiveWindow.ScrollRow = 16
  $excel.ActiveWindow.ScrollColumn = 2
} catch {}
